$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-range cells (B8/B9) first so their shared-string
# slots are reused in place of the old "07:45.../07:50..." values.
$ws.Range("B8").Value = "17:45 - 17:49"
$ws.Range("B9").Value = "17:50 - 17:54"

# Replace the trigger text cells (A6/A7) with the new (unmarked, plain)
# prayer texts.
$ws.Range("A6").Value = "1. Покаяние^`n^`nОтче наш, Отец Небесный,^`nпрошу Тебя простить^`nвсе мои прегрешения,^`nвольные или невольные.^`n^`nАминь!^`n^`n2. Новая Молитва^`n^`nОТЧЕ НАШ, ОТЕЦ НЕБЕСНЫЙ,^`nЯ ПРИНИМАЮ ВЕРУ ТВОЮ,^`nОНА ЕСТЬ МОЙ ПУТЬ.^`n^`nЯ ПРИНИМАЮ КАНОНЫ ТВОИ ВЕЧНЫЕ,^`nС ЛЮБОВЬЮ К ТЕБЕ И ДЕЛАМ ТВОИМ,^`nПОДТВЕРЖДАЯ СВОЕЙ ЖИЗНЬЮ^`nВЕРНОСТЬ ТЕБЕ.^`n^`nГОСПОДИ, ПРОШУ ДАТЬ МНЕ НАДЕЖДУ^`nНА СПАСЕНИЕ ДУШИ МОЕЙ,^`nИ ДАРОВАТЬ МУДРОСТЬ ТВОЮ^`nДЛЯ ЖИЗНИ МОЕЙ ЗДЕСЬ,^`nНА ПЛАНЕТЕ СВЯТАЯ РУСЬ И В ВЕЧНОСТИ.^`n^`nПУСТЬ СВЯТА БУДЕТ УВЕРЕННОСТЬ МОЯ,^`nЧТО ТЫ ЕСМЬ!^`n^`nГосподи, я Люблю Тебя, Благодарю Тебя и Уповаю на Милость Твою! Аминь!"

$ws.Range("A7").Value = "Отче наш, Отец Небесный! Волею Создателя, Пророка и Народа Пространство Святая Русь ЕСМЬ Равенство и Любовь Навечно! Да будет Свет Истины!"

# Move the window scroll position / active selection to match the
# author's final view (best-effort; scroll position itself is not
# exposed through this headless host, only the active selection is).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F6").Select() | Out-Null
